$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 now has an actual reported value instead of a forecasted one.
# Replace the forecast formula in I38 with the hardcoded actual figure,
# matching the styling used for the other "actual" rows (I30:I37).
$ws.Range("I38").Value = 311357
$ws.Range("I37").Copy()
$ws.Range("I38").PasteSpecial(-4122)

$wb.Application.Calculate()
